# Apply the edits described by the commit diff to the active workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Column width changes
#    col D (4): 93 -> 62 ; col M (13): 19 -> 13
#    (COM ColumnWidth is stored-width minus 5/6; use that offset so the
#    saved OOXML <col width=".."> lands on the exact target integer.)
# ---------------------------------------------------------------------
$ws.Columns("D:D").ColumnWidth = 61.1666666667
$ws.Columns("M:M").ColumnWidth = 12.1666666667

# ---------------------------------------------------------------------
# 2. Header block (rows 6-12)
# ---------------------------------------------------------------------
$ws.Range("N6").Value = 45232

$ws.Range("D7").Value = "TP WATERS"
$ws.Range("N7").Value = "IREN231102-1"

$ws.Range("D8").Value = "MONASTEREVIN RD, KILNAGORNAN, CO. KILDARE, R51 EV29, IRELAND"

$ws.Range("D9").Value = "(PH) 353 85 386 6717 / (Email) marie@tpwaters.ie" + [char]0x2019

$ws.Range("D10").Value = "SAME AS CONSIGNEE"

$ws.Range("D11").Value = ""
$ws.Range("D12").Value = ""

# ---------------------------------------------------------------------
# 3. Vehicle table (rows 17-21)
#    NOTE: column H (RECNO) values are purely-numeric-looking strings
#    that must stay stored as TEXT (as in the source file). Briefly force
#    a text number format so the literal isn't coerced into a number,
#    then clear the format again so no stray style index is left behind.
# ---------------------------------------------------------------------
function Set-TextValue($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# Row 17 - vehicle 1
$ws.Range("E17").Value = 2018
$ws.Range("F17").Value = "VOLKSWAGEN"
$ws.Range("G17").Value = "GOLF"
Set-TextValue $ws.Range("H17") "257434"
$ws.Range("I17").Value = "WVWZZZAUZJW276698"
$ws.Range("J17").Value = 1240
$ws.Range("K17").Value = 4.26
$ws.Range("L17").Value = 1.8
$ws.Range("M17").Value = 1.48
$ws.Range("N17").Value = 11.349
$ws.Range("O17").Value = "1200 CC"

# Row 18 - vehicle 2
$ws.Range("E18").Value = 2018
$ws.Range("F18").Value = "MAZDA"
$ws.Range("G18").Value = "DEMIO"
Set-TextValue $ws.Range("H18") "257667"
$ws.Range("I18").Value = "DJ3FS-518869"
$ws.Range("J18").Value = 1030
$ws.Range("K18").Value = 4.06
$ws.Range("L18").Value = 1.69
$ws.Range("M18").Value = 1.5
$ws.Range("N18").Value = 10.292
$ws.Range("O18").Value = "1300 CC"

# Row 19 - vehicle 3
$ws.Range("E19").Value = 2017
$ws.Range("F19").Value = "VOLKSWAGEN"
$ws.Range("G19").Value = "GOLF"
Set-TextValue $ws.Range("H19") "259249"
$ws.Range("I19").Value = "WVWZZZAUZJW045727"
$ws.Range("J19").Value = 1240
$ws.Range("K19").Value = 4.26
$ws.Range("L19").Value = 1.8
$ws.Range("M19").Value = 1.48
$ws.Range("N19").Value = 11.349
$ws.Range("O19").Value = "1200 CC"

# Row 20 - vehicle 4
$ws.Range("E20").Value = 2018
$ws.Range("F20").Value = "VOLKSWAGEN"
$ws.Range("G20").Value = "GOLF"
Set-TextValue $ws.Range("H20") "259997"
$ws.Range("I20").Value = "WVWZZZAUZJW296218"
$ws.Range("J20").Value = 1240
$ws.Range("K20").Value = 4.26
$ws.Range("L20").Value = 1.8
$ws.Range("M20").Value = 1.48
$ws.Range("N20").Value = 11.349
$ws.Range("O20").Value = "1200 CC"

# Row 21 - vehicle 5
$ws.Range("E21").Value = 2018
$ws.Range("F21").Value = "VOLKSWAGEN"
$ws.Range("G21").Value = "POLO"
Set-TextValue $ws.Range("H21") "262279"
$ws.Range("I21").Value = "WVWZZZAWZJU028693"
$ws.Range("J21").Value = 1160
$ws.Range("K21").Value = 4.06
$ws.Range("L21").Value = 1.75
$ws.Range("M21").Value = 1.45
$ws.Range("N21").Value = 10.302
$ws.Range("O21").Value = "990 CC"

# ---------------------------------------------------------------------
# 4. Totals row (22)
# ---------------------------------------------------------------------
$ws.Range("I22").Value = 5910
$ws.Range("M22").Value = 54.641

# ---------------------------------------------------------------------
# 5. Booking / shipping block (rows 24-29)
# ---------------------------------------------------------------------
$ws.Range("G24").Value = "EBKG07048809"

$ws.Range("D27").Value = "MSC NAGOYA V"

$ws.Range("D28").Value = "NO.HI346A"

$ws.Range("D29").Value = 45246
